$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text value would otherwise be auto-converted to a number
# by Excel's type inference are pre-formatted as Text so the literal string is kept.
$ws.Range('D2').Value = '29.617.57'
$ws.Range('E2').Value = '  +0.98%  '
$ws.Range('D3').Value = '1.851.15'
$ws.Range('E3').Value = '  +0.28%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9986'
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '240.70'
$ws.Range('E5').Value = '  +0.21%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.6307'
$ws.Range('E6').Value = '  +0.45%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.9999'
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.07475'
$ws.Range('E8').Value = '  -1.22%  '
$ws.Range('E9').Value = '  +0.38%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '25.13'
$ws.Range('E10').Value = '  +2.66%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07751'
$ws.Range('E11').Value = '  +0.04%  '
$ws.Range('D12').Value = '1.850.20'
$ws.Range('E12').Value = '  +0.24%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.032'
$ws.Range('E13').Value = '  +0.66%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.6824'
$ws.Range('E14').Value = '  +0.71%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.00001032'
$ws.Range('E15').Value = '  -0.64%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '82.81'
$ws.Range('E16').Value = '  -0.24%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '6.327'
$ws.Range('E17').Value = '  +3.83%  '
$ws.Range('D18').Value = '29.610.39'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '230.63'
$ws.Range('E19').Value = '  +0.63%  '
$ws.Range('E20').Value = '  +0.74%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.9996'
$ws.Range('E21').Value = '  +0.03%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '7.541'
$ws.Range('E22').Value = '  +1.59%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.000'
$ws.Range('E23').Value = '  +0.11%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '159.53'
$ws.Range('E24').Value = '  +0.39%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '8.522'
$ws.Range('E25').Value = '  +1.10%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.1368'
$ws.Range('E26').Value = '  -1.87%  '
$ws.Range('E27').Value = '  -0.29%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.06592'
$ws.Range('E28').Value = '  +16.11%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.464'
$ws.Range('E29').Value = '  +2.68%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.487'
$ws.Range('E30').Value = '  +1.08%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.114'
$ws.Range('E31').Value = '  +0.18%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.105'
$ws.Range('E32').Value = '  +1.59%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.853'
$ws.Range('E33').Value = '  +1.81%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.145'
$ws.Range('E34').Value = '  -0.72%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.6997'
$ws.Range('E35').Value = '  +0.41%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.567'
$ws.Range('E36').Value = '  -0.42%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.01868'
$ws.Range('E37').Value = '  +2.10%  '
$ws.Range('B38').Value = 'MXToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.841'
$ws.Range('E38').Value = '  +4.60%  '
$ws.Range('B39').Value = 'Maker'
$ws.Range('C39').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D39').Value = '1.258.70'
$ws.Range('E39').Value = '  +1.82%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '6.774'
$ws.Range('E40').Value = '  +5.65%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.9422'
$ws.Range('E41').Value = '  +4.81%  '
$ws.Range('E42').Value = '  +0.22%  '
$ws.Range('E43').Value = '  -0.05%  '
$ws.Range('E44').Value = '  +0.10%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '66.24'
$ws.Range('E45').Value = '  +1.11%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.737'
$ws.Range('E46').Value = '  +3.88%  '
$ws.Range('E47').Value = '  -0.13%  '
$ws.Range('B48').Value = 'Algorand'
$ws.Range('C48').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.1163'
$ws.Range('E48').Value = '  +1.13%  '
$ws.Range('B49').Value = 'BabyDogeCoin'
$ws.Range('C49').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.00000000116'
$ws.Range('E49').Value = '  +1.45%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '9.017'
$ws.Range('E50').Value = '  +0.42%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.3954'
$ws.Range('E51').Value = '  -0.93%  '
